# Add Portuguese (por) zone_user master data rows to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: globaladmin
$ws.Range("A4").Value = "por"
$ws.Range("B4").Value = "IST"
$ws.Range("C4").Value = "globaladmin"
$ws.Range("D4").Value = $true

# Row 5: service-account-mosip-resident-client
$ws.Range("A5").Value = "por"
$ws.Range("B5").Value = "IST"
$ws.Range("C5").Value = "service-account-mosip-resident-client"
$ws.Range("D5").Value = $true

# Row 7 filled before row 6 so shared-string entries land in the same
# order as the source workbook ("officer" before "ganesh").
# Row 7: officer
$ws.Range("A7").Value = "por"
$ws.Range("B7").Value = "IST"
$ws.Range("C7").Value = "officer"
$ws.Range("D7").Value = $true

# Row 6: ganesh
$ws.Range("A6").Value = "por"
$ws.Range("B6").Value = "IST"
$ws.Range("C6").Value = "ganesh"
$ws.Range("D6").Value = $true

# Keep the same boolean "TRUE"/"FALSE" display format used by the
# existing is_active column.
$ws.Range("D4:D7").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# Mirror the final cell selection recorded in the workbook.
$ws.Range("D11").Select() | Out-Null
